$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!D2 -- "Latest Handoff Date" column: update timestamp (text, not a real date)
$wsOverview.Range("D2").Value = "2016-04-13 09:04:06"

# zh-cn!E2 -- "Latest Handoff Datetime" column: update timestamp (text, not a real date)
$wsZhCn.Range("E2").Value = "2016-03-13 09:04:03"

# de-de!E2 -- "Latest Handoff Datetime" column: update timestamp (text, not a real date)
$wsDeDe.Range("E2").Value = "2016-03-13 09:04:06"
